$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Helper: write one "match" row (columns A..N) on a results sheet and
# the rolling-average formulas in columns O and P.
# -----------------------------------------------------------------
function Add-MatchRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [int]$A, [int]$B, [int]$C, [int]$D, [int]$E,
        [string]$F,
        [int]$I, [int]$J, [int]$K, [int]$L, [int]$M, [int]$N,
        [int]$AvgStart, [int]$AvgEnd
    )

    $ws = $wb.Worksheets.Item($SheetName)

    $ws.Range("A$Row").Value = $A
    $ws.Range("B$Row").Value = $B
    $ws.Range("C$Row").Value = $C
    $ws.Range("D$Row").Value = $D
    $ws.Range("E$Row").Value = $E
    $ws.Range("F$Row").Value = $F
    $ws.Range("G$Row").Formula = "=_xlfn.XLOOKUP(F$Row,Rankings!B:B,Rankings!A:A)"
    $ws.Range("H$Row").Formula = "=_xlfn.XLOOKUP(F$Row, Rankings!B:B, Rankings!C:C)"
    $ws.Range("I$Row").Value = $I
    $ws.Range("J$Row").Value = $J
    $ws.Range("K$Row").Value = $K
    $ws.Range("L$Row").Value = $L
    $ws.Range("M$Row").Value = $M
    $ws.Range("N$Row").Value = $N
    $ws.Range("O$Row").Formula = "=AVERAGE(L$AvgStart`:L$AvgEnd)"
    $ws.Range("P$Row").Formula = "=AVERAGE(M$AvgStart`:M$AvgEnd)"
}

# -----------------------------------------------------------------
# Round 3 results finished coming in (26-06 matches) -> round 4
# predictions added: one new "opponent" row per country sheet.
# -----------------------------------------------------------------

# Romania (sheet16) vs Slovakia
Add-MatchRow -SheetName "Romania" -Row 18 `
    -A 17 -B 0 -C 4 -D 5 -E 1 -F "Slovakia" `
    -I 4 -J 3 -K 3 -L 1 -M 1 -N 1 -AvgStart 13 -AvgEnd 17

# Ukraine (sheet17) vs Belgium
Add-MatchRow -SheetName "Ukraine" -Row 18 `
    -A 17 -B 1 -C 5 -D 3 -E 2 -F "Belgium" `
    -I 4 -J 5 -K 1 -L 0 -M 0 -N 1 -AvgStart 13 -AvgEnd 17

# Belgium (sheet18) vs Ukraine
Add-MatchRow -SheetName "Belgium" -Row 21 `
    -A 20 -B 0 -C 4 -D 5 -E 1 -F "Ukraine" `
    -I 5 -J 3 -K 2 -L 0 -M 0 -N 1 -AvgStart 16 -AvgEnd 20

# Slovakia (sheet19) vs Romania
Add-MatchRow -SheetName "Slovakia" -Row 18 `
    -A 17 -B 1 -C 4 -D 3 -E 3 -F "Romania" `
    -I 4 -J 5 -K 1 -L 1 -M 1 -N 1 -AvgStart 13 -AvgEnd 17

# Turkey (sheet22) vs Czechia
Add-MatchRow -SheetName "Turkey" -Row 18 `
    -A 17 -B 0 -C 4 -D 5 -E 1 -F "Czechia" `
    -I 3 -J 6 -K 1 -L 2 -M 1 -N 1 -AvgStart 13 -AvgEnd 17

# Georgia (sheet23) vs Portugal
Add-MatchRow -SheetName "Georgia" -Row 17 `
    -A 16 -B 1 -C 5 -D 3 -E 2 -F "Portugal" `
    -I 3 -J 6 -K 1 -L 2 -M 0 -N 1 -AvgStart 12 -AvgEnd 16

# Portugal (sheet24) vs Georgia
Add-MatchRow -SheetName "Portugal" -Row 23 `
    -A 22 -B 0 -C 3 -D 6 -E 1 -F "Georgia" `
    -I 5 -J 3 -K 2 -L 0 -M 2 -N 1 -AvgStart 18 -AvgEnd 22

# Czechia (sheet25) vs Turkey
Add-MatchRow -SheetName "Czechia" -Row 18 `
    -A 17 -B 1 -C 3 -D 6 -E 1 -F "Turkey" `
    -I 4 -J 5 -K 1 -L 1 -M 2 -N 1 -AvgStart 13 -AvgEnd 17

# -----------------------------------------------------------------
# Restore the per-sheet selections left behind after entering the
# data (mirrors where the author's cursor ended up on each tab).
# -----------------------------------------------------------------
$wb.Worksheets.Item("Romania").Range("O17:P18").Select()
$wb.Worksheets.Item("Ukraine").Range("O17:P18").Select()
$wb.Worksheets.Item("Belgium").Range("O20:P21").Select()
$wb.Worksheets.Item("Slovakia").Range("O17:P18").Select()
$wb.Worksheets.Item("Georgia").Range("O16:P17").Select()
$wb.Worksheets.Item("Portugal").Range("J24").Select()
$wb.Worksheets.Item("Czechia").Range("N20").Select()

# Turkey is the sheet the author finished on, so activate it last and
# leave the cursor on I20 - this also makes it the workbook's active tab.
$ws = $wb.Worksheets.Item("Turkey")
$ws.Activate()
$ws.Range("I20").Select()
